# Update weekly Fruta/Hortaliza price data (Femacal de La Calera - Coco)
# Applies the per-row corrections to columns D (Fecha), M (Volumen),
# N/O/P (Precio min/max/promedio) and S (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = 44445; M = 45; N = 20000; O = 20000; P = 20000; S = 1000 }
    3  = @{ D = 44305; M = 20; N = 22000; O = 22000; P = 22000; S = 1100 }
    4  = @{ D = 44291; M = 70; N = 25000; O = 25000; P = 25000; S = 1250 }
    5  = @{ D = 44382; M = 24 }
    6  = @{ D = 44292; M = 30; N = 25000; O = 25000; P = 25000; S = 1250 }
    7  = @{ D = 44400; M = 45 }
    8  = @{ D = 44413; M = 45 }
    9  = @{ D = 44406; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    10 = @{ D = 44389; M = 20; N = 20000; O = 20000; P = 20000; S = 1000 }
    11 = @{ D = 44307; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    13 = @{ D = 44403; M = 50; N = 20000; O = 20000; P = 20000; S = 1000 }
    15 = @{ D = 44377; M = 25 }
    17 = @{ D = 44300; N = 22000; O = 22000; P = 22000; S = 1100 }
    19 = @{ D = 44448; M = 30; N = 22000; O = 22000; P = 22000; S = 1100 }
    20 = @{ D = 44376; M = 38 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
